$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume strings are written back as text (matching original inlineStr cells)
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "245.27"
$ws.Cells.Item(2, 5).Value = "-0.05%"

$ws.Range("E3").NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "4.44%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.118"
$ws.Cells.Item(4, 5).Value = "-0.23%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.05588"
$ws.Cells.Item(5, 5).Value = "0.00%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "6.475"
$ws.Cells.Item(6, 5).Value = "-0.56%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.8170"
$ws.Cells.Item(7, 5).Value = "-0.02%"

$ws.Range("E9").NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "-0.96%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.06993"
$ws.Cells.Item(10, 5).Value = "0.58%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.02889"
$ws.Cells.Item(11, 5).Value = "1.50%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.09376"
$ws.Cells.Item(12, 5).Value = "-0.02%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.001509"
$ws.Cells.Item(13, 5).Value = "-0.42%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "One"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(14, 4).Value = "0.0005999"
$ws.Cells.Item(14, 5).Value = "0.78%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "TigerCash"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15, 4).Value = "0.006152"
$ws.Cells.Item(15, 5).Value = "1.03%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16, 4).Value = "3.648"
$ws.Cells.Item(16, 5).Value = "4.21%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = "GateToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17, 4).Value = "3.037"
$ws.Cells.Item(17, 5).Value = "0.65%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = "BTSEToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(18, 4).Value = "2.183"
$ws.Cells.Item(18, 5).Value = "5.82%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.03138"
$ws.Cells.Item(20, 5).Value = "-0.63%"

$ws.Range("E21").NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "-2.23%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.743"
$ws.Cells.Item(22, 5).Value = "-0.06%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04597"
$ws.Cells.Item(23, 5).Value = "-2.87%"

$ws.Range("E24").NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-0.13%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.001245"
$ws.Cells.Item(25, 5).Value = "-0.22%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.004514"
$ws.Cells.Item(26, 5).Value = "-2.86%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.00009598"
$ws.Cells.Item(27, 5).Value = "-1.08%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.03641"
$ws.Cells.Item(40, 5).Value = "-0.66%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 4).Value = "0.1361"
$ws.Cells.Item(41, 5).Value = "29.32%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = "KickToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(42, 4).Value = "0.006169"
$ws.Cells.Item(42, 5).Value = "-0.38%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.002660"
$ws.Cells.Item(43, 5).Value = "1.31%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.008872"
$ws.Cells.Item(44, 5).Value = "6.94%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00005349"
$ws.Cells.Item(45, 5).Value = "0.96%"

$ws.Range("E46").NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "-0.05%"

$ws.Range("E48").NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "20.34%"

$ws.Range("E49").NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "-0.05%"

$ws.Range("E50").NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "-0.05%"
